$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Coordinator Role paragraph: append two new runs at paragraph end.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(7)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter(" Does one third of the project and ensures the project ")

$p = $d.Paragraphs(7)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter("tasks and files are properly compiled.")

# ---------------------------------------------------------------------
# 2) Supervisor Role paragraph: split "send remainder" -> "send " / "the
#    remainder" / ", if necessary..." and append two more trailing runs.
# ---------------------------------------------------------------------
$hit = $d.Content
$found = $hit.Find.Execute("remainder, if necessary", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$delStart = $hit.Start
$p = $d.Paragraphs(8)
$r = $p.Range
$delEnd = $r.End - 1
$delRng = $d.Range($delStart, $delEnd)
$delRng.Delete()

$p = $d.Paragraphs(8)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter("the remainder")

$p = $d.Paragraphs(8)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter(", if necessary, with specifics on areas that need attention, and effective focus.")

$p = $d.Paragraphs(8)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter(" ")

$p = $d.Paragraphs(8)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter(". Does one third of the project and ensures the project tasks and files are properly compiled.")

# ---------------------------------------------------------------------
# 3) Insert a blank paragraph right after the Supervisor Role paragraph.
#    A placeholder char is typed immediately after the break (inheriting
#    plain formatting) and then removed, which leaves a truly empty
#    <w:p/> rather than one carrying the next paragraph's bold run.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(8)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter("X")

$p = $d.Paragraphs(8)
$r = $p.Range
$boundary = $r.End - 2
$brk = $d.Range($boundary, $boundary)
$brk.InsertParagraphAfter()

$blank = $d.Paragraphs(9)
$br = $blank.Range
$delX = $d.Range($br.Start, $br.Start + 1)
$delX.Delete()

# ---------------------------------------------------------------------
# 4) Project Submitter Role paragraph (now index 10): append three runs.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(10)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter(" ")

$p = $d.Paragraphs(10)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter("Does one third of the project and ensures the project tasks and files are properly compiled.")

$p = $d.Paragraphs(10)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter(" Finally, submits the project assignment as per the instructions provided.")

# ---------------------------------------------------------------------
# 5) Insert a blank paragraph right after the Project Submitter Role
#    paragraph, using the same placeholder technique as step 3.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(10)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter("X")

$p = $d.Paragraphs(10)
$r = $p.Range
$boundary = $r.End - 2
$brk = $d.Range($boundary, $boundary)
$brk.InsertParagraphAfter()

$blank2 = $d.Paragraphs(11)
$br2 = $blank2.Range
$delX2 = $d.Range($br2.Start, $br2.Start + 1)
$delX2.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
